# Atualizacao automatica: 2025-08-21 21:00:24
#
# The upstream data-collection job re-paginated the detections that sit in
# rows 7-11 of the active sheet: each row's detection (Fly_ID, image,
# placa, location, lat/long, bounding box and confidence) moved up into
# the row above it, with what used to be row 7 wrapping around into row
# 11. Columns B (Class) and C (First_Detection_Date) are identical for
# every row in this block, so they do not need to be touched.
#
# This script reproduces that rotation using the existing workbook data
# (it reads the current values out of the block and writes them back one
# row higher, wrapping the top row's values around to the bottom), rather
# than hard-coding the values, so it is resilient to being re-applied.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 7
$lastRow = 11

# Columns whose contents actually vary row to row within this block.
$textCols = @("A", "D", "E", "F")          # plain text fields
$numericCols = @("G", "H")                 # latitude / longitude (numbers)
$stringNumberCols = @("I", "J")            # bounding box + confidence -
                                            # stored as *text* even though
                                            # they look numeric, so they
                                            # must not be auto-converted.

# ---- 1. Snapshot the current ("before") values for the whole block -------
$before = @{}
foreach ($col in $textCols + $numericCols + $stringNumberCols) {
    $colValues = @{}
    for ($r = $firstRow; $r -le $lastRow; $r++) {
        $colValues[$r] = $ws.Range($col + $r).Value2
    }
    $before[$col] = $colValues
}

# ---- 2. Helper that returns the source row for a rotated destination -----
function Get-SourceRow([int]$destRow) {
    if ($destRow -lt $lastRow) {
        return $destRow + 1
    } else {
        return $firstRow
    }
}

# ---- 3. Write the rotated values back -------------------------------------
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = Get-SourceRow $r

    foreach ($col in $textCols) {
        $ws.Range($col + $r).Value = $before[$col][$srcRow]
    }

    foreach ($col in $numericCols) {
        $ws.Range($col + $r).Value2 = $before[$col][$srcRow]
    }

    foreach ($col in $stringNumberCols) {
        $cell = $ws.Range($col + $r)
        # Force text storage so strings such as "0.62" or "702,633,740,690"
        # are not reinterpreted as numbers by Excel's normal "smart" entry
        # parsing, then drop the temporary number format override so the
        # cell's style is left exactly as it was (default/general).
        $cell.NumberFormat = "@"
        $cell.Value = $before[$col][$srcRow]
        $cell.ClearFormats()
    }
}
